# Update "Baseline model performance matrix_SAM.xlsx":
#  - add a "Team Member" column (F) marking who ran each row (Sam / Arjun)
#  - widen column B (augmentation text) since new rows have longer text
#  - add new VGG16 / extra SimpleCNN_SAM result rows contributed by Arjun
#  - a couple of the copied-down rows are plain formulas referencing the row above

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- widen column B to fit the longer augmentation description ----
$ws.Columns.Item(2).ColumnWidth = 73.34

# ---- new "Team Member" header + best-fit-ish width for column F ----
$ws.Cells.Item(1, 6).Value = "Team Member"
$ws.Cells.Item(1, 6).Font.Bold = $true
$ws.Columns.Item(6).ColumnWidth = 12.42

# ---- tag the existing 9 rows (all Sam's work) with the new column ----
$ws.Cells.Item(2, 6).Value = "Sam"
$ws.Cells.Item(3, 6).Value = "Sam"
$ws.Cells.Item(4, 6).Value = "Sam"
$ws.Cells.Item(5, 6).Value = "Sam"
$ws.Cells.Item(6, 6).Value = "Sam"
$ws.Cells.Item(7, 6).Value = "Sam"
$ws.Cells.Item(8, 6).Value = "Sam"
$ws.Cells.Item(9, 6).Value = "Sam"
$ws.Cells.Item(10, 6).Value = "Sam"

# rows 5-9 shrink a touch (85 -> 80), row 10 collapses a lot (85 -> 32)
$ws.Rows.Item(5).RowHeight = 80
$ws.Rows.Item(6).RowHeight = 80
$ws.Rows.Item(7).RowHeight = 80
$ws.Rows.Item(8).RowHeight = 80
$ws.Rows.Item(9).RowHeight = 80
$ws.Rows.Item(10).RowHeight = 32

# ---- row 11: extra SimpleCNN_SAM run (Arjun) ----
$ws.Cells.Item(11, 1).Value = "SimpleCNN_SAM"
$ws.Cells.Item(11, 2).Value = "Yes(Resize, RandomResizeCrop,`nHorizontal and vertical flip,Rotation and Gaussian Blur)"
$ws.Cells.Item(11, 3).Value = 0.66381987577599999
$ws.Cells.Item(11, 4).Value = 0.97699999999999998
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = "Arjun"
# this row keeps the sheet's default (un-wrapped) row height
$ws.Rows.Item(11).AutoFit()

# ---- row 12: filled down from row 11 via formulas ----
$ws.Cells.Item(12, 1).Formula = "=A11"
$ws.Cells.Item(12, 2).Formula = "=B11"
$ws.Cells.Item(12, 3).Value = 0.66381987570000001
$ws.Cells.Item(12, 4).Value = 0.85199999999999998
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = "Arjun"

# ---- rows 13-16: new VGG16 results (Arjun) ----
$ws.Cells.Item(13, 1).Value = "VGG16"
$ws.Cells.Item(13, 2).Value = "Yes(Resize, RandomResizeCrop,`nHorizontal and vertical flip,Rotation and Gaussian Blur)  Requires_grad=False"
$ws.Cells.Item(13, 2).WrapText = $true
$ws.Cells.Item(13, 3).Value = 0.66537000000000002
$ws.Cells.Item(13, 4).Value = 0.94
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = "Arjun"
$ws.Rows.Item(13).RowHeight = 32

$ws.Cells.Item(14, 1).Value = "VGG16"
$ws.Cells.Item(14, 2).Value = "Yes(Resize, RandomResizeCrop,`nHorizontal and vertical flip,Rotation and Gaussian Blur)  Requires_grad=False"
$ws.Cells.Item(14, 2).WrapText = $true
$ws.Cells.Item(14, 3).Value = 0.71273291925465798
$ws.Cells.Item(14, 4).Value = 0.86329871416091897
$ws.Cells.Item(14, 5).Value = 5
$ws.Cells.Item(14, 6).Value = "Arjun"
$ws.Rows.Item(14).RowHeight = 32

$ws.Cells.Item(15, 1).Formula = "=A14"
$ws.Cells.Item(15, 2).Value = "Yes(Resize, RandomResizeCrop,`nHorizontal and vertical flip,Rotation and Gaussian Blur)  Requires_grad=False"
$ws.Cells.Item(15, 2).WrapText = $true
$ws.Cells.Item(15, 3).Value = 0.78649068322981297
$ws.Cells.Item(15, 4).Value = 0.76341283321380604
$ws.Cells.Item(15, 5).Value = 20
$ws.Cells.Item(15, 6).Formula = "=F14"
$ws.Rows.Item(15).RowHeight = 32

$ws.Cells.Item(16, 1).Formula = "=A15"
$ws.Cells.Item(16, 2).Value = "Yes(Resize, RandomResizeCrop,`nHorizontal and vertical flip,Rotation and Gaussian Blur)  Requires_grad=False"
$ws.Cells.Item(16, 2).WrapText = $true
$ws.Cells.Item(16, 3).Value = 0.75
$ws.Cells.Item(16, 4).Value = 0.70205193758010798
$ws.Cells.Item(16, 5).Value = 40
$ws.Cells.Item(16, 6).Value = "Arjun"
$ws.Rows.Item(16).RowHeight = 32

# ---- leave the view roughly where the author left it ----
$ws.Range("C24").Select()
